$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "UNcwb833"
$ws.Range("C2").Value = "ohncjxw43"
$ws.Range("D2").Value = "b7Gw5`$R&"
$ws.Range("F2").Value = "gqadzRme"
$ws.Range("G2").Value = "vmcO"
$ws.Range("B2").Value = 23111405
